$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Felipe" row (row 5) entirely
$ws.Rows.Item(5).Delete()

# Delete the two oldest week columns (B=11_02_2024, C=18_02_2024) - data shifts left
$ws.Range("B1:C1").EntireColumn.Delete()

# After deletion, the new last-used column is E (was F). Add new column E = 17_03_2024
$ws.Range("E1").Value = "17_03_2024"
$ws.Range("E2").Value = 1131
$ws.Range("E3").Value = 1070
$ws.Range("E4").Value = 1569
$ws.Range("E5").Value = 180

# Update selection to match the target view
$ws.Range("E6").Select()
